$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column price cells keep their original text formatting
# (e.g. "7.63", "1.40", "0.0000171") instead of being auto-converted
# to numbers by Excel when assigned via .Value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.991.64"
$ws.Range("E2").Value = "  -0.65%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.390.33"
$ws.Range("E3").Value = "  -1.49%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.22"
$ws.Range("E5").Value = "  -0.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.14"
$ws.Range("E6").Value = "  -1.54%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.390.68"
$ws.Range("E7").Value = "  -1.51%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.475"
$ws.Range("E9").Value = "  -0.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.63"
$ws.Range("E10").Value = "  +0.21%  "

$ws.Range("E11").Value = "  -2.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.395"
$ws.Range("E12").Value = "  +1.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.971.64"
$ws.Range("E13").Value = "  -1.45%  "

$ws.Range("E14").Value = "  +2.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.82"
$ws.Range("E15").Value = "  -3.25%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.397.97"
$ws.Range("E16").Value = "  -1.54%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000171"
$ws.Range("E17").Value = "  -1.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.013.40"
$ws.Range("E18").Value = "  -0.81%  "

$ws.Range("E19").Value = "  -3.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.79"
$ws.Range("E20").Value = "  -4.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.91"
$ws.Range("E21").Value = "  -4.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "382.46"
$ws.Range("E22").Value = "  -4.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.556"
$ws.Range("E23").Value = "  -1.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.42"
$ws.Range("E24").Value = "  +0.79%  "

$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000117"
$ws.Range("E26").Value = "  -4.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.534.72"
$ws.Range("E27").Value = "  -1.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.181"
$ws.Range("E28").Value = "  +0.93%  "

$ws.Range("E29").Value = "  -0.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.35"
$ws.Range("E30").Value = "  -2.86%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.97"
$ws.Range("E31").Value = "  -3.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.15"
$ws.Range("E32").Value = "  -1.43%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.40"
$ws.Range("E33").Value = "  -6.23%  "

$ws.Range("E34").Value = "  -0.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.43"
$ws.Range("E35").Value = "  -2.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.98"
$ws.Range("E36").Value = "  -0.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "167.10"
$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.423.75"
$ws.Range("E38").Value = "  -1.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.98"
$ws.Range("E39").Value = "  -3.00%  "

$ws.Range("E40").Value = "  -4.61%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0768"
$ws.Range("E41").Value = "  -2.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.28"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.781"
$ws.Range("E43").Value = "  -2.85%  "

$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("E45").Value = "  -2.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.67"
$ws.Range("E46").Value = "  -3.61%  "

$ws.Range("E47").Value = "  -1.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.483.66"
$ws.Range("E48").Value = "  -4.94%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.80"
$ws.Range("E49").Value = "  -2.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.88"
$ws.Range("E50").Value = "  -0.85%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0263"
$ws.Range("E51").Value = "  +0.61%  "
